$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Text-level edits: turn
#      " attribute that will contain the subsequent visit, if there is no
#        subsequent visit (this is the first visit) it will contain "None""
#    into
#      " attribute that will contain the previous visit that lead to this
#        visit, if there is no previous visit (this is the first visit) it
#        will contain "None""
# ------------------------------------------------------------------

# Replace the first occurrence of "subsequent" -> "previous"
$r1 = $d.Content
$r1.Find.Execute("subsequent", $true, $false, $false, $false, $false, $true, 1, $false, "previous", 1) | Out-Null

# Replace the (new) second occurrence of "subsequent" -> "previous"
$r2 = $d.Content
$r2.Find.Execute("subsequent", $true, $false, $false, $false, $false, $true, 1, $false, "previous", 1) | Out-Null

# Insert " that lead to this visit" right after the first "previous visit" (before the comma)
$r3 = $d.Content
$r3.Find.Execute("previous visit,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertionPoint = $d.Range($r3.End - 1, $r3.End - 1)
$insertionPoint.InsertBefore(" that lead to this visit")

# ------------------------------------------------------------------
# 2. Re-split the now-merged run into the seven runs shown in the diff by
#    nudging (and restoring) a formatting property across each sub-range;
#    this makes the engine materialize separate <w:r> elements at those
#    boundaries while leaving the effective formatting unchanged.
# ------------------------------------------------------------------

$locator = $d.Content
$locator.Find.Execute(" attribute that will contain the previous visit that lead to this visit, if there is no previous visit (this is the first visit) it will contain", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $locator.Start

$lq = [char]0x201C
$rq = [char]0x201D

$segments = @(
    " attribute that will contain the ",
    "previous",
    " visit",
    " that lead to this visit",
    ", if there is no ",
    "previous",
    (" visit (this is the first visit) it will contain " + $lq + "None" + $rq)
)

$pos = $base
foreach ($seg in $segments) {
    $segLen = $seg.Length
    $segRange = $d.Range($pos, $pos + $segLen)
    $segRange.Font.Bold = $true
    $segRange.Font.Bold = $false
    $pos += $segLen
}

Write-Output "Edit complete"
